# 977: Add GS to extract process and GS tab to example files
#
# Adds a new "GS" worksheet (sheetId 9) after the existing last sheet
# ("CMS"), with a header row of 9 columns reusing existing shared strings
# where possible, styled with a new Calibri/12/black font, and makes it
# the active/selected sheet (so CMS loses tabSelected).

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after the current last sheet (CMS),
# so it lands at the end of the tab strip and becomes the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$gs = $wb.Worksheets.Add($null, $lastSheet)
$gs.Name = "GS"

# Header row for the new GS sheet.
$headers = @(
    "Contact_ID",
    "Contact_Date",
    "Contact_Type_Code",
    "Contact_Type_Desc",
    "OM_Name",
    "OM_Key",
    "OM_Grade",
    "OM_Team_Key",
    "OM_Provider_Code"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $gs.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Style + select the header range; this both marks GS as the active sheet
# (tabSelected) and gives the header cells their black Calibri font.
$headerRange = $gs.Range("A1:I1")
$headerRange.Font.Color = 0
$headerRange.Select()
